$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to hold a literal text value (matches the source
    # workbook's inline-string "Price" column) even when the text looks
    # like a number. The leading apostrophe makes Excel store it as text;
    # resetting the style back to Normal avoids leaving behind the
    # "number stored as text" formatting Excel would otherwise apply.
    $cell = $ws.Range($range)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Column D price refreshes
Set-TextValue "D2"  "236.09"
Set-TextValue "D3"  "21.99"
Set-TextValue "D4"  "5.361"
Set-TextValue "D5"  "0.05608"
Set-TextValue "D6"  "6.474"
Set-TextValue "D7"  "3.342"
Set-TextValue "D8"  "0.7990"
Set-TextValue "D9"  "1.046"
Set-TextValue "D11" "0.07300"
Set-TextValue "D12" "0.03148"
Set-TextValue "D13" "0.02967"
Set-TextValue "D14" "0.09242"
Set-TextValue "D15" "0.001668"
Set-TextValue "D16" "3.256"
Set-TextValue "D17" "0.04786"

Set-TextValue "D18" "0.0005714"
$ws.Range("E18").Value = "17OneONEWorstin24h"

Set-TextValue "D19" "0.006217"
Set-TextValue "D20" "0.005065"
Set-TextValue "D21" "0.001051"
Set-TextValue "D22" "0.0001502"
Set-TextValue "D23" "0.0003701"

Set-TextValue "D24" "3.968"
$ws.Range("E24").Value = "23LEOLEOBestin24h"

Set-TextValue "D25" "2.203"

Set-TextValue "D40" "0.04088"
Set-TextValue "D41" "0.007033"

# Row 42 / row 43 swap (CEJI <-> BKEXToken trade places in the ranking)
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1038"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003003"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue "D44" "0.008831"
Set-TextValue "D45" "0.00005437"

Set-TextValue "D47" "0.6756"

Set-TextValue "D48" "0.03662"
$ws.Range("E48").Value = "47BOLOBOLO"

Set-TextValue "D50" "0.01011"
